# The scraper originally only grabbed team statistics; this update adds the
# season record (Wins / Losses / Ties) columns AC:AE to the roster table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new "Wins"/"Losses"/"Ties" columns, matching the
#     existing header formatting (bold, bordered, centered) used by AB1. ---
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# --- Data rows (2-45): every player on this roster shares the team's
#     1991 season record: 98 wins, 64 losses, 0 ties. ---
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 29).Value = 98   # AC
    $ws.Cells.Item($row, 30).Value = 64   # AD
    $ws.Cells.Item($row, 31).Value = 0    # AE
}
